$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2..59, columns B|C|D|F (E is always blank).
# Column A already contains a simple 0-based index (row-2) in the source
# workbook and that pattern is unchanged by this edit, except it now needs
# to extend down to row 59.
$rowData = @"
NSE:AIRAN|NSE:HEIDELBERG|NSE:ABFRL|NSE:COROMANDEL
NSE:ALPHAETF|NSE:MUFIN||NSE:MOTHERSON
NSE:ANANDRATHI|NSE:PANSARI||
NSE:ANIKINDS|||
NSE:ASHOKA|||
NSE:ATL|||
NSE:AXISCADES|||
NSE:BAJAJELEC|||
NSE:BALMLAWRIE|||
NSE:BFSI|||
NSE:CANFINHOME|||
NSE:CARBORUNIV|||
NSE:CERA|||
NSE:COCHINSHIP|||
NSE:CONSUMBEES|||
NSE:COROMANDEL|||
NSE:DATAPATTNS|||
NSE:ESG|||
NSE:FACT|||
NSE:GIPCL|||
NSE:GMRINFRA|||
NSE:GPIL|||
NSE:GPPL|||
NSE:GRSE|||
NSE:GULFPETRO|||
NSE:HAL|||
NSE:HDFCBANK|||
NSE:HDFCMOMENT|||
NSE:HDFCNEXT50|||
NSE:HERCULES|||
NSE:HFCL|||
NSE:IDEAFORGE|||
NSE:IGPL|||
NSE:IMPAL|||
NSE:JSWENERGY|||
NSE:KAPSTON|||
NSE:KAYNES|||
NSE:KBCGLOBAL|||
NSE:KICL|||
NSE:M&M|||
NSE:MIDHANI|||
NSE:MOMENTUM|||
NSE:MOMOMENTUM|||
NSE:MONIFTY500|||
NSE:MOTHERSON|||
NSE:MTARTECH|||
NSE:NACLIND|||
NSE:NDLVENTURE|||
NSE:NETF|||
NSE:NFL|||
NSE:NIF100BEES|||
NSE:NIFTYQLITY|||
NSE:PHOENIXLTD|||
NSE:PILANIINVS|||
NSE:PNBHOUSING|||
NSE:PRAKASH|||
NSE:QUICKHEAL|||
NSE:RADIANTCMS|||
"@

$lines = $rowData -split "`n"

# First, extend column A (the 0-based index column) down through row 59,
# copying the existing style (bold, centered, thin-bordered) from A42 so the
# new rows 43:59 look like the existing ones.
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A43:A59").PasteSpecial(-4122) | Out-Null

$r = 2
foreach ($line in $lines) {
    $parts = $line -split '\|', -1
    $bVal = $parts[0]
    $cVal = $parts[1]
    $dVal = $parts[2]
    $fVal = $parts[3]

    $ws.Cells.Item($r, 1).Value = ($r - 2)
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = $fVal

    $r = $r + 1
}

Write-Host "Done. Last row written:" ($r - 1)
